$d = $word.ActiveDocument

$d.Content.Find.Execute("2025-08-06 Wednesday", $true, $false, $false, $false, $false, $true, 1, $false, "2025-08-07 Thursday", 2) | Out-Null
$d.Content.Find.Execute("367×8=", $true, $false, $false, $false, $false, $true, 1, $false, "699×4=", 2) | Out-Null
$d.Content.Find.Execute("552×7=", $true, $false, $false, $false, $false, $true, 1, $false, "820×7=", 2) | Out-Null
$d.Content.Find.Execute("795×4=", $true, $false, $false, $false, $false, $true, 1, $false, "803×6=", 2) | Out-Null
$d.Content.Find.Execute("578×4=", $true, $false, $false, $false, $false, $true, 1, $false, "650×8=", 2) | Out-Null
$d.Content.Find.Execute("206×3=", $true, $false, $false, $false, $false, $true, 1, $false, "285×6=", 2) | Out-Null
$d.Content.Find.Execute("953×4=", $true, $false, $false, $false, $false, $true, 1, $false, "624×8=", 2) | Out-Null
$d.Content.Find.Execute("176×2=", $true, $false, $false, $false, $false, $true, 1, $false, "610×9=", 2) | Out-Null
$d.Content.Find.Execute("629×5=", $true, $false, $false, $false, $false, $true, 1, $false, "849×5=", 2) | Out-Null
$d.Content.Find.Execute("898×2=", $true, $false, $false, $false, $false, $true, 1, $false, "873×7=", 2) | Out-Null
$d.Content.Find.Execute("437×4=", $true, $false, $false, $false, $false, $true, 1, $false, "491×9=", 2) | Out-Null
$d.Content.Find.Execute("884×8=", $true, $false, $false, $false, $false, $true, 1, $false, "321×8=", 2) | Out-Null
$d.Content.Find.Execute("348×8=", $true, $false, $false, $false, $false, $true, 1, $false, "412×6=", 2) | Out-Null
$d.Content.Find.Execute("525×5=", $true, $false, $false, $false, $false, $true, 1, $false, "536×7=", 2) | Out-Null
$d.Content.Find.Execute("390×9=", $true, $false, $false, $false, $false, $true, 1, $false, "993×2=", 2) | Out-Null
$d.Content.Find.Execute("994×5=", $true, $false, $false, $false, $false, $true, 1, $false, "840×4=", 2) | Out-Null
$d.Content.Find.Execute("539×6=", $true, $false, $false, $false, $false, $true, 1, $false, "737×6=", 2) | Out-Null
$d.Content.Find.Execute("643×5=", $true, $false, $false, $false, $false, $true, 1, $false, "465×4=", 2) | Out-Null
$d.Content.Find.Execute("235×4=", $true, $false, $false, $false, $false, $true, 1, $false, "154×8=", 2) | Out-Null
$d.Content.Find.Execute("448×3=", $true, $false, $false, $false, $false, $true, 1, $false, "950×4=", 2) | Out-Null
$d.Content.Find.Execute("515×4=", $true, $false, $false, $false, $false, $true, 1, $false, "408×4=", 2) | Out-Null
$d.Content.Find.Execute("338×6=", $true, $false, $false, $false, $false, $true, 1, $false, "775×3=", 2) | Out-Null
$d.Content.Find.Execute("495×8=", $true, $false, $false, $false, $false, $true, 1, $false, "481×6=", 2) | Out-Null
$d.Content.Find.Execute("879×9=", $true, $false, $false, $false, $false, $true, 1, $false, "720×3=", 2) | Out-Null
$d.Content.Find.Execute("229×4=", $true, $false, $false, $false, $false, $true, 1, $false, "526×5=", 2) | Out-Null
$d.Content.Find.Execute("918×3=", $true, $false, $false, $false, $false, $true, 1, $false, "861×2=", 2) | Out-Null
